# Scheduled market-data refresh: update cached Universalis price snapshots
# (currentAveragePrice / NQ / HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ columns)
# across the per-job Sheets. Generated by the commit-bot; values below are
# the new snapshot pulled for this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 102.4375
$ws.Range("I2").Value = 103.666664
$ws.Range("J2").Value = 98.75
$ws.Range("K2").Value = 103.666664
$ws.Range("L2").Value = 98.75
$ws.Range("M2").Value = 9.333336000000003
$ws.Range("N2").Value = -324.75
$ws.Range("H9").Value = 142.2
$ws.Range("I9").Value = 175.75
$ws.Range("K9").Value = 175.75
$ws.Range("M9").Value = -6.75
$ws.Range("H64").Value = 6111.1113
$ws.Range("I64").Value = 3333.3333
$ws.Range("K64").Value = 3333.3333
$ws.Range("M64").Value = -3085.3333
$ws.Range("H67").Value = 6111.1113
$ws.Range("I67").Value = 3333.3333
$ws.Range("K67").Value = 3333.3333
$ws.Range("M67").Value = -2475.3333
$ws.Range("H111").Value = 1089.2
$ws.Range("I111").Value = 933.3333
$ws.Range("J111").Value = 1323
$ws.Range("K111").Value = 2799.9999
$ws.Range("L111").Value = 3969
$ws.Range("M111").Value = 267.0001000000002
$ws.Range("N111").Value = -10103
$ws.Range("H116").Value = 3996
$ws.Range("I116").Value = 3995
$ws.Range("J116").Value = 3996.6667
$ws.Range("K116").Value = 3995
$ws.Range("L116").Value = 3996.6667
$ws.Range("M116").Value = -553
$ws.Range("N116").Value = -10880.6667
$ws.Range("H132").Value = 11934.632
$ws.Range("I132").Value = 11934.632
$ws.Range("K132").Value = 35803.896
$ws.Range("M132").Value = -33273.896
$ws.Range("H135").Value = 867.4286
$ws.Range("I135").Value = 677.7692
$ws.Range("K135").Value = 6099.922799999999
$ws.Range("M135").Value = -3564.922799999999
$ws.Range("H137").Value = 2272.3333
$ws.Range("I137").Value = 1354.75
$ws.Range("K137").Value = 4064.25
$ws.Range("M137").Value = -1514.25
$ws.Range("H138").Value = 8627.857
$ws.Range("I138").Value = 348.5
$ws.Range("J138").Value = 11939.6
$ws.Range("K138").Value = 1045.5
$ws.Range("L138").Value = 35818.8
$ws.Range("M138").Value = 4094.5
$ws.Range("N138").Value = -46098.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1333.1666
$ws.Range("J88").Value = 1602.75
$ws.Range("L88").Value = 1602.75
$ws.Range("N88").Value = -2414.75
$ws.Range("H91").Value = 1333.1666
$ws.Range("J91").Value = 1602.75
$ws.Range("L91").Value = 1602.75
$ws.Range("N91").Value = -4410.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4563.6665
$ws.Range("I54").Value = 3395.75
$ws.Range("J54").Value = 6899.5
$ws.Range("K54").Value = 3395.75
$ws.Range("L54").Value = 6899.5
$ws.Range("M54").Value = -2911.75
$ws.Range("N54").Value = -7867.5
$ws.Range("H61").Value = 40000
$ws.Range("J61").Value = 40000
$ws.Range("L61").Value = 40000
$ws.Range("N61").Value = -40626
$ws.Range("H64").Value = 767.2
$ws.Range("J64").Value = 795.3333
$ws.Range("L64").Value = 795.3333
$ws.Range("N64").Value = -1245.3333
$ws.Range("H67").Value = 767.2
$ws.Range("J67").Value = 795.3333
$ws.Range("L67").Value = 795.3333
$ws.Range("N67").Value = -2355.3333
$ws.Range("H95").Value = 16874.334
$ws.Range("J95").Value = 16874.334
$ws.Range("L95").Value = 16874.334
$ws.Range("N95").Value = -22366.334
$ws.Range("H134").Value = 3573.75
$ws.Range("I134").Value = 1048.75
$ws.Range("K134").Value = 3146.25
$ws.Range("M134").Value = -611.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1444.8235
$ws.Range("I58").Value = 839.10345
$ws.Range("K58").Value = 839.10345
$ws.Range("M58").Value = -636.10345
$ws.Range("H99").Value = 4980.1665
$ws.Range("I99").Value = 4981.6665
$ws.Range("K99").Value = 4981.6665
$ws.Range("M99").Value = -3483.6665
$ws.Range("H126").Value = 4980.1665
$ws.Range("I126").Value = 4981.6665
$ws.Range("K126").Value = 14944.9995
$ws.Range("M126").Value = -12474.9995
$ws.Range("H134").Value = 3390.762
$ws.Range("I134").Value = 2624.7778
$ws.Range("K134").Value = 7874.3334
$ws.Range("M134").Value = -5339.3334
$ws.Range("H136").Value = 1444.8235
$ws.Range("I136").Value = 839.10345
$ws.Range("K136").Value = 2517.31035
$ws.Range("M136").Value = 32.68965000000026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3576.6
$ws.Range("I132").Value = 1900
$ws.Range("K132").Value = 17100
$ws.Range("M132").Value = -14570
$ws.Range("H140").Value = 2183.3076
$ws.Range("J140").Value = 4472.1665
$ws.Range("L140").Value = 13416.4995
$ws.Range("N140").Value = -23776.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1500
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H122").Value = 3532
$ws.Range("I122").Value = 3292.3333
$ws.Range("K122").Value = 9876.999899999999
$ws.Range("M122").Value = -7426.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5052
$ws.Range("I40").Value = 5057.579
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 5057.579
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -4921.579
$ws.Range("N40").Value = -5271
$ws.Range("H46").Value = 5603.2
$ws.Range("I46").Value = 724.5
$ws.Range("J46").Value = 6353.769
$ws.Range("K46").Value = 724.5
$ws.Range("L46").Value = 6353.769
$ws.Range("M46").Value = -536.5
$ws.Range("N46").Value = -6729.769
$ws.Range("H54").Value = 26666.666
$ws.Range("J54").Value = 26666.666
$ws.Range("L54").Value = 26666.666
$ws.Range("N54").Value = -27954.666
$ws.Range("H100").Value = 9400
$ws.Range("I100").Value = 7000
$ws.Range("K100").Value = 7000
$ws.Range("M100").Value = -6459
$ws.Range("H132").Value = 7137.9443
$ws.Range("I132").Value = 7335.533
$ws.Range("K132").Value = 22006.599
$ws.Range("M132").Value = -19476.599
$ws.Range("H136").Value = 893.2
$ws.Range("I136").Value = 893.2
$ws.Range("K136").Value = 2679.6
$ws.Range("M136").Value = -129.6000000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10284.214
$ws.Range("I18").Value = 9998.25
$ws.Range("J18").Value = 12000
$ws.Range("K18").Value = 9998.25
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = -9825.25
$ws.Range("N18").Value = -12346
$ws.Range("H96").Value = 1310.6666
$ws.Range("I96").Value = 1362
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 1362
$ws.Range("L96").Value = 900
$ws.Range("M96").Value = 11
$ws.Range("N96").Value = -3646
$ws.Range("H97").Value = 13333.333
$ws.Range("J97").Value = 13333.333
$ws.Range("L97").Value = 13333.333
$ws.Range("N97").Value = -15315.333
$ws.Range("H122").Value = 5237.6665
$ws.Range("I122").Value = 4167.8
$ws.Range("K122").Value = 12503.4
$ws.Range("M122").Value = -10053.4
$ws.Range("H136").Value = 3469.647
$ws.Range("I136").Value = 2248.375
$ws.Range("J136").Value = 4555.222
$ws.Range("K136").Value = 6745.125
$ws.Range("L136").Value = 13665.666
$ws.Range("M136").Value = -4195.125
$ws.Range("N136").Value = -18765.666
